# Feature to copy Java Requirements list and append to each associate
#
# The sheet held a ".NET" cohort list (ID + first-name pairs in columns A/B,
# rows 2-6, with row 1 / rows 7-8 reserved but blank). This edit repurposes
# row 1 as a header ("Week  1" list item) and refreshes the body rows with
# the new "Java Requirements" checklist entries (Diahandra / Doing / Done),
# while keeping the last two existing rows (Zaur, Diahandra) intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: previously only A1 existed (blank placeholder). Give it real
# content and add a new B1 cell alongside it.
$ws.Range("A1").Value = "62b757f7a5d42e8dd3afd466"
$ws.Range("B1").Value = "Week  1"

# B1 is a brand-new cell - copy the formatting already used by the rest of
# column B (e.g. B2) onto it so it carries the same cell style.
$ws.Range("B2").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats

# Rows 2-4: replace with the new Java Requirements rows.
$ws.Range("A2").Value = "62b7576afe12f938a57c624b"
$ws.Range("B2").Value = "Diahandra"

$ws.Range("A3").Value = "62b666a6bc1aa85d93c63b7f"
$ws.Range("B3").Value = "Doing"

$ws.Range("A4").Value = "62b666a6bc1aa85d93c63b80"
$ws.Range("B4").Value = "Done"

# Rows 5-6: unchanged content (Zaur, Diahandra) - rewritten for parity.
$ws.Range("A5").Value = "62b74c21e17fdb80e8513e7a"
$ws.Range("B5").Value = "Zaur"

$ws.Range("A6").Value = "62b74c23e2197787f1b7e3b4"
$ws.Range("B6").Value = "Diahandra"
